# Auto-generated edit script: update cryptocurrency price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.417.07'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '1.900.51'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'325.26"
$ws.Range("E5").Value = '  -2.78%  '
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").Value = "'0.4816"
$ws.Range("E7").Value = '  +3.19%  '
$ws.Range("D8").Value = "'0.4068"
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").Value = "'0.08061"
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("D10").Value = "'1.002"
$ws.Range("D11").Value = "'23.30"
$ws.Range("E11").Value = '  +4.29%  '
$ws.Range("D12").Value = '1.929.06'
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("D13").Value = "'5.942"
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = "'7.060"
$ws.Range("E14").Value = '  -1.57%  '
$ws.Range("D15").Value = "'89.77"
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = "'0.06708"
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = "'0.00001031"
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = "'17.56"
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").Value = '29.415.16'
$ws.Range("E21").Value = '  -0.65%  '
$ws.Range("D22").Value = "'5.532"
$ws.Range("E22").Value = '  -0.88%  '
$ws.Range("D23").Value = "'11.77"
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("D24").Value = "'2.158"
$ws.Range("E24").Value = '  -2.57%  '
$ws.Range("D25").Value = '2.098.33'
$ws.Range("E25").Value = '  -4.43%  '
$ws.Range("D26").Value = "'155.28"
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("E27").Value = '  -0.43%  '
$ws.Range("D28").Value = "'6.054"
$ws.Range("D29").Value = "'2.092"
$ws.Range("E29").Value = '  -2.09%  '
$ws.Range("D30").Value = "'118.46"
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("E31").Value = '  -3.30%  '
$ws.Range("D32").Value = "'0.09510"
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").Value = "'3.542"
$ws.Range("E33").Value = '  -1.08%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = "'1.388"
$ws.Range("E34").Value = '  -3.00%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = "'5.407"
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").Value = "'0.02249"
$ws.Range("E36").Value = '  -0.86%  '
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("D38").Value = "'1.175"
$ws.Range("E38").Value = '  -0.74%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = "'0.5866"
$ws.Range("E39").Value = '  -0.61%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = "'7.872"
$ws.Range("E40").Value = '  -6.56%  '
$ws.Range("D41").Value = "'0.1841"
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").Value = "'10.22"
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = "'2.402"
$ws.Range("E43").Value = '  +1.72%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = "'1.280"
$ws.Range("E44").Value = '  +1.44%  '
$ws.Range("D45").Value = "'0.07752"
$ws.Range("E45").Value = '  +3.14%  '
$ws.Range("D46").Value = "'12.23"
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").Value = "'0.5519"
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("D49").Value = "'113.37"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = "'72.06"
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = "'0.2930"
$ws.Range("E51").Value = '  -2.07%  '
